$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.83711850643158
$ws.Range("B1").Value = 2.274204969406128
$ws.Range("C1").Value = 2.286669254302979
$ws.Range("D1").Value = 2.639342546463013
$ws.Range("E1").Value = 3.424020767211914
